$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Card")
$ws.Activate()

# Adding more property classes: Go and Community Chest rows
# Row 30: Go -- acts as its own "Type", Property Value 200, Location 0
$ws.Range("A30").Value = "Go"
$ws.Range("B30").Value = "Go"
$ws.Range("C30").Value = 200
$ws.Range("D30").Value = 0

# Row 31: Community Chest
$ws.Range("A31").Value = "Community Chest"

# Match the author's final view/selection state
$excel.ActiveWindow.ScrollRow = 19
$ws.Range("A32").Select()
